# Emissions ceramics + update notebook residential tertiary
$wb = $excel.ActiveWorkbook

# --- Sheet "0D": drop the obsolete "init_installed_capacity" row and refresh
#     the init_energy_need_per_ton value ---
$ws0D = $wb.Worksheets.Item("0D")
$ws0D.Rows.Item(3).Delete()
$ws0D.Range("B2").Value = 0.27088069255578107

# --- Sheet "Production_system": rename the capacity column and refresh its
#     values for the updated production systems ---
$wsProd = $wb.Worksheets.Item("Production_system")
$wsProd.Range("B1").Value = "init_unite_prod"
$wsProd.Range("B2").Value = 16530.162694442846
$wsProd.Range("B5").Value = 1394.1101067602399
$wsProd.Range("B6").Value = 181.05326061821299

# --- Update selections / active sheet to match the saved UI state ---
$ws0D.Range("D7").Select()

$wsProd.Activate()
$wsProd.Range("C30").Select()
